# Turkey super-lig 2023-2024: re-sync a handful of match rows that were
# previously out of chronological/page order on betexplorer, and append
# the Istanbulspor AS vs Galatasaray fixture that was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several row-pairs had their match data (columns F..V: home/away teams,
# scores, odds, timestamps, url) swapped relative to each other. Columns
# A..E (Indice, pais, torneio, temporada, data_partida) are unaffected,
# so just exchange F..V between the two rows of each pair.
function Swap-MatchData($ws, $row1, $row2) {
    for ($col = 6; $col -le 22; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $val1 = $cell1.Value2
        $val2 = $cell2.Value2
        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

Swap-MatchData $ws 10 11
Swap-MatchData $ws 15 16
Swap-MatchData $ws 29 30
Swap-MatchData $ws 49 50

# Append the new match row (Istanbulspor AS 0 - 1 Galatasaray) as row 58.
$newRow = 58

$ws.Cells.Item($newRow, 1).Value = 57
$ws.Cells.Item($newRow, 2).Value = "turkey"
$ws.Cells.Item($newRow, 3).Value = "super-lig"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45195.79166666666
$ws.Cells.Item($newRow, 6).Value = "Istanbulspor AS"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Galatasaray"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 5.63
$ws.Cells.Item($newRow, 11).Value = "20/08/2023 00:12"
$ws.Cells.Item($newRow, 12).Value = 10
$ws.Cells.Item($newRow, 13).Value = "26/09/2023 18:59"
$ws.Cells.Item($newRow, 14).Value = 4.68
$ws.Cells.Item($newRow, 15).Value = "20/08/2023 00:12"
$ws.Cells.Item($newRow, 16).Value = 6.12
$ws.Cells.Item($newRow, 17).Value = "26/09/2023 18:59"
$ws.Cells.Item($newRow, 18).Value = 1.55
$ws.Cells.Item($newRow, 19).Value = "20/08/2023 00:12"
$ws.Cells.Item($newRow, 20).Value = 1.29
$ws.Cells.Item($newRow, 21).Value = "26/09/2023 18:57"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/turkey/super-lig/istanbulspor-as-galatasaray/6PyFWoHk/"

# Match the formatting of the existing data rows: column A (Indice) is
# bold/bordered/centered, column E (data_partida) is a date-time number
# format. Copy those styles from the row directly above onto the new row.
$ws.Range("A57").Copy() | Out-Null
$ws.Range("A58").PasteSpecial(-4122) | Out-Null

$ws.Range("E57").Copy() | Out-Null
$ws.Range("E58").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
